$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numeric-looking strings (e.g. "42.744.76",
# "7.70"). Mark it as Text before writing so Excel keeps the literal digits
# (trailing zeros, multi-dot "thousands" groupings) instead of auto-coercing to a number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '42.744.76'
$ws.Cells.Item(2, 5).Value = '  -6.69%  '
$ws.Cells.Item(3, 4).Value = '2.545.94'
$ws.Cells.Item(3, 5).Value = '  -5.13%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).Value = '299.15'
$ws.Cells.Item(5, 5).Value = '  -4.20%  '
$ws.Cells.Item(6, 4).Value = '94.12'
$ws.Cells.Item(6, 5).Value = '  -4.62%  '
$ws.Cells.Item(7, 4).Value = '0.573'
$ws.Cells.Item(7, 5).Value = '  -4.26%  '
$ws.Cells.Item(8, 5).Value = '  +0.09%  '
$ws.Cells.Item(9, 4).Value = '0.547'
$ws.Cells.Item(9, 5).Value = '  -6.29%  '
$ws.Cells.Item(10, 4).Value = '36.07'
$ws.Cells.Item(10, 5).Value = '  -6.11%  '
$ws.Cells.Item(11, 4).Value = '0.0804'
$ws.Cells.Item(11, 5).Value = '  -5.21%  '
$ws.Cells.Item(12, 4).Value = '7.70'
$ws.Cells.Item(12, 5).Value = '  -5.75%  '
$ws.Cells.Item(13, 4).Value = '0.115'
$ws.Cells.Item(13, 5).Value = '  +6.39%  '
$ws.Cells.Item(14, 4).Value = '2.935.70'
$ws.Cells.Item(14, 5).Value = '  -4.94%  '
$ws.Cells.Item(15, 4).Value = '2.562.03'
$ws.Cells.Item(15, 5).Value = '  -4.23%  '
$ws.Cells.Item(16, 4).Value = '0.877'
$ws.Cells.Item(16, 5).Value = '  -6.16%  '
$ws.Cells.Item(17, 4).Value = '14.20'
$ws.Cells.Item(17, 5).Value = '  -6.74%  '
$ws.Cells.Item(18, 4).Value = '42.748.12'
$ws.Cells.Item(18, 5).Value = '  -6.65%  '
$ws.Cells.Item(19, 4).Value = '12.71'
$ws.Cells.Item(19, 5).Value = '  -1.55%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0980'
$ws.Cells.Item(20, 5).Value = '  -4.07%  '
$ws.Cells.Item(21, 4).Value = '6.55'
$ws.Cells.Item(21, 5).Value = '  -4.75%  '
$ws.Cells.Item(22, 4).Value = '71.68'
$ws.Cells.Item(22, 5).Value = '  -4.87%  '
$ws.Cells.Item(23, 4).Value = '252.98'
$ws.Cells.Item(23, 5).Value = '  -10.74%  '
$ws.Cells.Item(24, 4).Value = '2.91'
$ws.Cells.Item(24, 5).Value = '  -5.07%  '
$ws.Cells.Item(25, 5).Value = '  -5.01%  '
$ws.Cells.Item(26, 4).Value = '28.95'
$ws.Cells.Item(26, 5).Value = '  -7.96%  '
$ws.Cells.Item(27, 5).Value = '  +0.15%  '
$ws.Cells.Item(28, 4).Value = '10.14'
$ws.Cells.Item(28, 5).Value = '  -4.26%  '
$ws.Cells.Item(29, 4).Value = '36.92'
$ws.Cells.Item(29, 5).Value = '  -3.74%  '
$ws.Cells.Item(30, 5).Value = '  -5.31%  '
$ws.Cells.Item(31, 4).Value = '6.05'
$ws.Cells.Item(31, 5).Value = '  -3.55%  '
$ws.Cells.Item(32, 4).Value = '152.19'
$ws.Cells.Item(32, 5).Value = '  -1.99%  '
$ws.Cells.Item(33, 5).Value = '  -2.45%  '
$ws.Cells.Item(34, 4).Value = '2.15'
$ws.Cells.Item(34, 5).Value = '  -9.41%  '
$ws.Cells.Item(35, 4).Value = '3.37'
$ws.Cells.Item(35, 5).Value = '  -10.66%  '
$ws.Cells.Item(36, 4).Value = '0.0793'
$ws.Cells.Item(36, 5).Value = '  -6.01%  '
$ws.Cells.Item(37, 5).Value = '  -5.91%  '
$ws.Cells.Item(38, 4).Value = '17.13'
$ws.Cells.Item(38, 5).Value = '  +5.26%  '
$ws.Cells.Item(39, 5).Value = '  -4.40%  '
$ws.Cells.Item(40, 4).Value = '23.11'
$ws.Cells.Item(40, 5).Value = '  -11.32%  '
$ws.Cells.Item(41, 2).Value = 'VeChain'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(41, 4).Value = '0.0310'
$ws.Cells.Item(41, 5).Value = '  -5.58%  '
$ws.Cells.Item(42, 2).Value = 'NEARProtocol'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(42, 4).Value = '3.40'
$ws.Cells.Item(42, 5).Value = '  -6.09%  '
$ws.Cells.Item(43, 4).Value = '3.86'
$ws.Cells.Item(43, 5).Value = '  -3.14%  '
$ws.Cells.Item(44, 4).Value = '2.081.23'
$ws.Cells.Item(44, 5).Value = '  -3.40%  '
$ws.Cells.Item(45, 4).Value = '0.998'
$ws.Cells.Item(45, 5).Value = '  +0.17%  '
$ws.Cells.Item(46, 5).Value = '  +3.61%  '
$ws.Cells.Item(47, 4).Value = '9.04'
$ws.Cells.Item(47, 5).Value = '  -3.61%  '
$ws.Cells.Item(48, 4).Value = '84.32'
$ws.Cells.Item(48, 5).Value = '  -10.75%  '
$ws.Cells.Item(49, 4).Value = '105.30'
$ws.Cells.Item(49, 5).Value = '  -6.29%  '
$ws.Cells.Item(50, 4).Value = '2.791.23'
$ws.Cells.Item(50, 5).Value = '  -5.02%  '
$ws.Cells.Item(51, 5).Value = '  -3.32%  '

Write-Output 'Updated cryptos list'
